$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Add()
$ws.Name = "TestSheet"
$ws.Range("A1").Value = "Hello"
